$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '58.872.31'
$ws.Range("E2").Value = '  +4.38%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.304.91'
$ws.Range("E3").Value = '  +2.45%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.16%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '400.72'
$ws.Range("E5").Value = '  +0.65%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '109.98'
$ws.Range("E6").Value = '  -0.99%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.588'
$ws.Range("E7").Value = '  +6.09%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  -0.01%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.638'
$ws.Range("E9").Value = '  +2.88%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.88'
$ws.Range("E10").Value = '  +1.40%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0984'
$ws.Range("E11").Value = '  +7.07%  '

# Row 12
$ws.Range("E12").Value = '  +1.51%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.818.47'
$ws.Range("E13").Value = '  +2.18%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '8.42'
$ws.Range("E14").Value = '  +3.91%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '19.31'
$ws.Range("E15").Value = '  +1.35%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.293.42'
$ws.Range("E16").Value = '  +2.13%  '

# Row 17
$ws.Range("E17").Value = '  -0.12%  '

# Row 18
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '58.504.29'
$ws.Range("E18").Value = '  +4.06%  '

# Row 19
$ws.Range("B19").Value = 'Uniswap'
$ws.Range("C19").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.89'
$ws.Range("E19").Value = '  -0.14%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.32'
$ws.Range("E20").Value = '  -0.55%  '

# Row 21
$ws.Range("E21").Value = '  +7.06%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '13.01'
$ws.Range("E22").Value = '  +0.08%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '304.21'
$ws.Range("E23").Value = '  +1.92%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '74.75'
$ws.Range("E24").Value = '  -1.13%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.19'
$ws.Range("E25").Value = '  -1.00%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '28.41'
$ws.Range("E26").Value = '  +1.14%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '4.43'
$ws.Range("E27").Value = '  +1.28%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.86'
$ws.Range("E28").Value = '  -3.44%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.39'
$ws.Range("E29").Value = '  -0.64%  '

# Row 30
$ws.Range("E30").Value = '  -1.23%  '

# Row 31
$ws.Range("B31").Value = 'Hedera'
$ws.Range("C31").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.114'
$ws.Range("E31").Value = '  +3.31%  '

# Row 32
$ws.Range("B32").Value = 'Dai'
$ws.Range("C32").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.00'
$ws.Range("E32").Value = '  -0.44%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.43'
$ws.Range("E33").Value = '  +2.51%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '40.06'
$ws.Range("E34").Value = '  +8.98%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0529'
$ws.Range("E35").Value = '  +7.22%  '

# Row 36
$ws.Range("E36").Value = '  +0.03%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '51.91'
$ws.Range("E37").Value = '  +1.05%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.36'
$ws.Range("E38").Value = '  +7.99%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.999'
$ws.Range("E39").Value = '  -0.02%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.49'
$ws.Range("E40").Value = '  -1.03%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '137.65'
$ws.Range("E41").Value = '  +2.12%  '

# Row 42
$ws.Range("E42").Value = '  +2.60%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.89'
$ws.Range("E43").Value = '  -1.76%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.94'
$ws.Range("E44").Value = '  -1.43%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '16.81'
$ws.Range("E45").Value = '  -3.82%  '

# Row 46
$ws.Range("E46").Value = '  -1.78%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.31'
$ws.Range("E47").Value = '  +11.03%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '22.54'
$ws.Range("E48").Value = '  +1.24%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.166.17'
$ws.Range("E49").Value = '  +1.67%  '

# Row 50
$ws.Range("E50").Value = '  -0.12%  '

# Row 51
$ws.Range("E51").Value = '  -13.33%  '
